$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (changed) date column (C) for rows 2-10
# from serial date 45233 (2023-11-03) to 45243 (2023-11-13)
for ($row = 2; $row -le 10; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45233) {
        $cell.Value = 45243
    }
}
